$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Efna3","Epha4","ECs",1,0.3333333333333333,0.046374,0.139122,0.6592866045237633,0.6592866045237632,3,1,6.030956000000001,18.092868,0.364814105361131,0.3648141053611309,0.279679553544,2.517115981896,0.2405170528059145,0.2405170528059144),
  @("ECs","Efna3","Epha4","FAPs",1,0.3333333333333333,0.046374,0.139122,0.6592866045237633,0.6592866045237632,3,1,6.789877333333333,20.369632,0.4107214552505144,0.4107214552505143,0.314873771456,2.833863943104,0.2707831536371704,0.2707831536371703),
  @("ECs","Efna3","Epha4","MuSCs",1,0.3333333333333333,0.046374,0.139122,0.6592866045237633,0.6592866045237632,3,1,3.710753333333333,11.13226,0.2244644393883547,0.2244644393883547,0.17208247508,1.54874227572,0.1479863980806785,0.1479863980806784),
  @("MuSCs","Efna3","Epha4","ECs",1,0.3333333333333333,0.02396566666666667,0.071897,0.3407133954762367,0.3407133954762367,3,1,6.030956000000001,18.092868,0.364814105361131,0.3648141053611309,0.1445358811773333,1.300822930596,0.1242970525552165,0.1242970525552165),
  @("MuSCs","Efna3","Epha4","FAPs",1,0.3333333333333333,0.02396566666666667,0.071897,0.3407133954762367,0.3407133954762367,3,1,6.789877333333333,20.369632,0.4107214552505144,0.4107214552505143,0.1627239368782222,1.464515431904,0.139938301613344,0.1399383016133439),
  @("MuSCs","Efna3","Epha4","MuSCs",1,0.3333333333333333,0.02396566666666667,0.071897,0.3407133954762367,0.3407133954762367,3,1,3.710753333333333,11.13226,0.2244644393883547,0.2244644393883547,0.08893067746888887,0.8003760972199999,0.07647804130767627,0.07647804130767626)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

Write-Host "done"
